$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- values from original row 49
$ws.Range("D2").Value = 44321
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 3333

# Row 3 <- values from original row 14
$ws.Range("D3").Value = 44412
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10500
$ws.Range("M3").Value = 10260
$ws.Range("P3").Value = 3420

# Row 4 <- values from original row 11
$ws.Range("D4").Value = 44503
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("P4").Value = 2833

# Row 5 <- values from original row 20
$ws.Range("D5").Value = 44293
$ws.Range("J5").Value = 16
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 3333

# Row 6 <- values from original row 4
$ws.Range("D6").Value = 44419
$ws.Range("J6").Value = 16
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 3333

# Row 7 <- values from original row 38
$ws.Range("D7").Value = 44266
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("P7").Value = 3333

# Row 8 <- values from original row 3
$ws.Range("D8").Value = 44517
$ws.Range("J8").Value = 16
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 9500
$ws.Range("P8").Value = 3167

# Row 9 <- values from original row 37
$ws.Range("D9").Value = 44447
$ws.Range("J9").Value = 16
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10500
$ws.Range("M9").Value = 10250
$ws.Range("P9").Value = 3417

# Row 10 <- values from original row 47
$ws.Range("D10").Value = 44468
$ws.Range("J10").Value = 16
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 10500
$ws.Range("P10").Value = 3500

# Row 11 <- values from original row 16
$ws.Range("D11").Value = 44559
$ws.Range("J11").Value = 7
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 11143
$ws.Range("P11").Value = 3714

# Row 12 <- values from original row 35
$ws.Range("D12").Value = 44545
$ws.Range("J12").Value = 25
$ws.Range("K12").Value = 9000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 9480
$ws.Range("P12").Value = 3160

# Row 13 <- values from original row 2
$ws.Range("D13").Value = 44342
$ws.Range("J13").Value = 17
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 3333

# Row 14 <- values from original row 28
$ws.Range("D14").Value = 44580
$ws.Range("J14").Value = 16
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9500
$ws.Range("P14").Value = 3167

# Row 15 <- values from original row 39
$ws.Range("D15").Value = 44391
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 10000
$ws.Range("P15").Value = 3333

# Row 16 <- values from original row 9
$ws.Range("D16").Value = 44307
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("P16").Value = 3333

# Row 17 <- values from original row 10
$ws.Range("D17").Value = 44181
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 11000
$ws.Range("P17").Value = 3667

# Row 18 <- values from original row 22
$ws.Range("D18").Value = 44405
$ws.Range("J18").Value = 16
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10500
$ws.Range("M18").Value = 10250
$ws.Range("P18").Value = 3417

# Row 19 <- values from original row 5
$ws.Range("D19").Value = 44377
$ws.Range("J19").Value = 16
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10500
$ws.Range("M19").Value = 10250
$ws.Range("P19").Value = 3417

# Row 20 <- values from original row 24
$ws.Range("D20").Value = 44587
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9500
$ws.Range("P20").Value = 3167

# Row 21 <- values from original row 6
$ws.Range("D21").Value = 44384
$ws.Range("J21").Value = 25
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10500
$ws.Range("M21").Value = 10260
$ws.Range("P21").Value = 3420

# Row 22 <- values from original row 21
$ws.Range("D22").Value = 44435
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 10500
$ws.Range("M22").Value = 10250
$ws.Range("P22").Value = 3417

# Row 23 <- values from original row 46
$ws.Range("D23").Value = 44314
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 10000
$ws.Range("P23").Value = 3333

# Row 24 <- values from original row 29
$ws.Range("D24").Value = 44510
$ws.Range("J24").Value = 16
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 9500
$ws.Range("P24").Value = 3167

# Row 25 <- values from original row 33
$ws.Range("D25").Value = 44540
$ws.Range("J25").Value = 32
$ws.Range("K25").Value = 8500
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 8719
$ws.Range("P25").Value = 2906

# Row 26 <- values from original row 19
$ws.Range("D26").Value = 44398
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 10500
$ws.Range("M26").Value = 10250
$ws.Range("P26").Value = 3417

# Row 27 <- values from original row 41
$ws.Range("D27").Value = 44328
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = 10000
$ws.Range("P27").Value = 3333

# Row 28 <- values from original row 23
$ws.Range("D28").Value = 44433
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 10500
$ws.Range("M28").Value = 10250
$ws.Range("P28").Value = 3417

# Row 29 <- values from original row 7
$ws.Range("D29").Value = 44363
$ws.Range("J29").Value = 16
$ws.Range("K29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 10000
$ws.Range("P29").Value = 3333

# Row 30 <- values from original row 27
$ws.Range("D30").Value = 44461
$ws.Range("J30").Value = 16
$ws.Range("K30").Value = 9500
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = 9750
$ws.Range("P30").Value = 3250

# Row 31 <- values from original row 44
$ws.Range("D31").Value = 44475
$ws.Range("J31").Value = 16
$ws.Range("K31").Value = 9000
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = 9500
$ws.Range("P31").Value = 3167

# Row 32 <- values from original row 26
$ws.Range("D32").Value = 44356
$ws.Range("J32").Value = 16
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("P32").Value = 3333

# Row 33 <- values from original row 12
$ws.Range("D33").Value = 44539
$ws.Range("J33").Value = 16
$ws.Range("K33").Value = 9000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 9500
$ws.Range("P33").Value = 3167

# Row 34 <- values from original row 13
$ws.Range("D34").Value = 44482
$ws.Range("J34").Value = 16
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 9500
$ws.Range("P34").Value = 3167

# Row 35 <- values from original row 34
$ws.Range("D35").Value = 44175
$ws.Range("J35").Value = 70
$ws.Range("K35").Value = 12000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 12000
$ws.Range("P35").Value = 4000

# Row 36 <- values from original row 8
$ws.Range("D36").Value = 44566
$ws.Range("J36").Value = 16
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 11000
$ws.Range("P36").Value = 3667

# Row 37 <- values from original row 15
$ws.Range("D37").Value = 44195
$ws.Range("J37").Value = 30
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = 10000
$ws.Range("P37").Value = 3333

# Row 38 <- values from original row 45
$ws.Range("D38").Value = 44300
$ws.Range("J38").Value = 16
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 10000
$ws.Range("P38").Value = 3333

# Row 39 <- values from original row 40
$ws.Range("D39").Value = 44454
$ws.Range("J39").Value = 16
$ws.Range("K39").Value = 9500
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 9750
$ws.Range("P39").Value = 3250

# Row 40 <- values from original row 36
$ws.Range("D40").Value = 44573
$ws.Range("J40").Value = 16
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 11000
$ws.Range("P40").Value = 3667

# Row 41 <- values from original row 32
$ws.Range("D41").Value = 44426
$ws.Range("J41").Value = 16
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 10500
$ws.Range("M41").Value = 10250
$ws.Range("P41").Value = 3417

# Row 42 <- values from original row 25
$ws.Range("D42").Value = 44524
$ws.Range("J42").Value = 16
$ws.Range("K42").Value = 9000
$ws.Range("L42").Value = 10000
$ws.Range("M42").Value = 9500
$ws.Range("P42").Value = 3167

# Row 43 <- values from original row 42
$ws.Range("D43").Value = 44489
$ws.Range("J43").Value = 16
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = 9500
$ws.Range("P43").Value = 3167

# Row 44 <- values from original row 30
$ws.Range("D44").Value = 44279
$ws.Range("J44").Value = 16
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = 10000
$ws.Range("P44").Value = 3333

# Row 45 <- values from original row 17
$ws.Range("D45").Value = 44370
$ws.Range("J45").Value = 16
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 10500
$ws.Range("M45").Value = 10250
$ws.Range("P45").Value = 3417

# Row 46 <- values from original row 48
$ws.Range("D46").Value = 44272
$ws.Range("J46").Value = 70
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = 10000
$ws.Range("P46").Value = 3333

# Row 47 <- values from original row 18
$ws.Range("D47").Value = 44349
$ws.Range("J47").Value = 12
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("M47").Value = 10000
$ws.Range("P47").Value = 3333

# Row 48 <- values from original row 43
$ws.Range("D48").Value = 44335
$ws.Range("J48").Value = 16
$ws.Range("K48").Value = 10000
$ws.Range("L48").Value = 10000
$ws.Range("M48").Value = 10000
$ws.Range("P48").Value = 3333

# Row 49 <- values from original row 31
$ws.Range("D49").Value = 44552
$ws.Range("J49").Value = 8
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 10000
$ws.Range("M49").Value = 10000
$ws.Range("P49").Value = 3333

